$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("A6:N6").ClearContents()
